$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-17 06:53:16"
$wsZh.Range("H2").Value = "2016-03-17 06:53:56"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-17 06:53:24"
$wsDe.Range("H2").Value = "2016-03-17 06:54:11"
